$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The target paragraph (about the dice) is a single big run in the source.
# We need to turn it into three runs:
#   "Nel gioco sono inoltre presenti dei dadi"
#   " (almeno uno)"                                            <- new text
#   ", ognuno con un valore minimo e un valore massimo (il valore
#      massimo non deve essere superiore a 6). "
# while leaving the rest of the paragraph's runs (and their text) untouched.
#
# This engine's text-mutation primitives (Find replace / InsertAfter /
# Range.Text =) normalize the whole host paragraph into a single run as a
# side effect, so after inserting the new phrase we re-establish every
# run boundary in the paragraph (old ones included) by toggling Bold on
# each exact sub-range -- a formatting-only edit, which this engine splits
# cleanly on Range boundaries without touching surrounding runs.
# ---------------------------------------------------------------------------

# Step 1: insert " (almeno uno)" right after "dadi" and before the comma.
$anchor = $d.Content
$anchor.Find.Execute("dei dadi, ognuno", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
if (-not $anchor.Find.Found) {
    throw "Could not locate the dice sentence to edit"
}
$insertAt = $anchor.Start + "dei dadi".Length
$insertionPoint = $d.Range($insertAt, $insertAt)
$insertionPoint.InsertAfter(" (almeno uno)")

# Step 2: re-find the paragraph (its start offset is unchanged by the
# insertion above, since the insertion happened after this anchor text).
$paraAnchor = $d.Content
$paraAnchor.Find.Execute("Nel gioco sono inoltre presenti", $true, $false, $false, `
                          $false, $false, $true, 1, $false, "", 0)
if (-not $paraAnchor.Find.Found) {
    throw "Could not re-locate the dice paragraph"
}
$paraStart = $paraAnchor.Start

# Step 3: the full, final text of the paragraph, split into the runs we
# want it to end up with (first three are new/changed, the rest restore
# the original run layout that got flattened by the edit above).
$segments = @(
    "Nel gioco sono inoltre presenti dei dadi",
    " (almeno uno)",
    ", ognuno con un valore minimo e un valore massimo (il valore massimo non deve essere superiore a 6). ",
    "Ogni squadra inizia il gioco con un certo numero di dadi a disposizione; tale ",
    "numero",
    " può essere modificato durante il corso della partita",
    " (può anche essere zero)",
    ", in particolare il punteggio che si ottiene rispondendo ai quiz/al task posti sulle caselle può modificare il numero di dadi a disposizione della squadra. ",
    "In",
    " tutte le caselle ",
    "di tipo normale",
    " ",
    "è previsto",
    " il lancio dei dadi."
)

# Turn the segment lengths into cumulative character offsets from the
# paragraph start: 0, len(seg0), len(seg0)+len(seg1), ...
$offsets = New-Object System.Collections.Generic.List[int]
$offsets.Add(0) | Out-Null
$running = 0
foreach ($seg in $segments) {
    $running = $running + $seg.Length
    $offsets.Add($running) | Out-Null
}

# Step 4: re-cut the run boundaries by toggling Bold on/off across each
# exact segment span -- this forces a run split at both ends of the span
# without altering any character or leaving residual formatting behind.
for ($i = 0; $i -lt $segments.Length; $i++) {
    $segStart = $paraStart + $offsets[$i]
    $segEnd = $paraStart + $offsets[$i + 1]
    $segRange = $d.Range($segStart, $segEnd)
    $segRange.Bold = 1
    $segRange.Bold = 0
}
